$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '28.385.53'
$ws.Range('E2').Value = '  +0.99%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '1.867.96'
$ws.Range('E3').Value = '  -0.36%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '1.020'
$ws.Range('E4').Value = '  +1.48%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '315.89'
$ws.Range('E5').Value = '  +0.89%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '1.017'
$ws.Range('E6').Value = '  +1.40%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '0.5118'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = '0.3922'
$ws.Range('E8').Value = '  +0.45%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.08283'
$ws.Range('E9').Value = '  -0.41%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '1.111'
$ws.Range('E10').Value = '  -0.70%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '41.91'
$ws.Range('E11').Value = '  +1.04%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '6.230'
$ws.Range('E12').Value = '  +0.43%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.875.63'
$ws.Range('E13').Value = '  -0.11%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = '20.41'
$ws.Range('E14').Value = '  -1.08%  '

$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').Value = '1.020'
$ws.Range('E15').Value = '  +1.32%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '7.204'
$ws.Range('E16').Value = '  -0.60%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001104'
$ws.Range('E17').Value = '  +0.37%  '

$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = '91.03'
$ws.Range('E18').Value = '  +0.09%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.06744'
$ws.Range('E19').Value = '  +1.39%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '17.63'
$ws.Range('E20').Value = '  -0.51%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.018'
$ws.Range('E21').Value = '  +1.36%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.957'
$ws.Range('E22').Value = '  -1.15%  '

$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').Value = '28.406.82'
$ws.Range('E23').Value = '  +0.91%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '11.11'
$ws.Range('E24').Value = '  -0.06%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.251'
$ws.Range('E25').Value = '  -0.29%  '

$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.087.85'
$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '161.00'
$ws.Range('E27').Value = '  +1.11%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '20.73'
$ws.Range('E28').Value = '  +0.71%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.365'
$ws.Range('E29').Value = '  -4.88%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '127.06'
$ws.Range('E30').Value = '  +1.75%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1050'
$ws.Range('E31').Value = '  -1.24%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.031'
$ws.Range('E32').Value = '  -0.55%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.794'
$ws.Range('E33').Value = '  -0.92%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.626'
$ws.Range('E34').Value = '  +0.42%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.02428'
$ws.Range('E35').Value = '  -0.86%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.06486'
$ws.Range('E36').Value = '  -1.08%  '

$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = '9.122'
$ws.Range('E37').Value = '  -5.08%  '

$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2174'
$ws.Range('E38').Value = '  -0.36%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.246'
$ws.Range('E39').Value = '  +1.53%  '

$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.183'
$ws.Range('E40').Value = '  -1.36%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6411'
$ws.Range('E41').Value = '  -1.02%  '

$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').Value = '4.976'
$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '11.13'
$ws.Range('E43').Value = '  -1.28%  '

$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '0.6013'
$ws.Range('E44').Value = '  -1.90%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '12.92'
$ws.Range('E45').Value = '  -1.00%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.689'
$ws.Range('E46').Value = '  +0.67%  '

$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.214'
$ws.Range('E47').Value = '  -5.47%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '1.985'
$ws.Range('E48').Value = '  -1.57%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.206'
$ws.Range('E49').Value = '  -2.01%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '121.54'
$ws.Range('E50').Value = '  +0.93%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06856'
$ws.Range('E51').Value = '  -0.68%  '

